$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 124
$ws1.Range("F3").Value = 1653
$ws1.Range("F6").Value = 428
$ws1.Range("F9").Value = 536

$ws4.Range("F2").Value = 124
$ws4.Range("F3").Value = 1653
$ws4.Range("F6").Value = 428
$ws4.Range("F9").Value = 536
